$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Short Term"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Short Term")

# Update existing rows 121-126
$ws1.Range("B121").Value = -9.19
$ws1.Range("C121").Value = -9.61
$ws1.Range("D121").Value = -1.55
$ws1.Range("F121").Value = 14.18
$ws1.Range("G121").Value = -0.87

$ws1.Range("B122").Value = 39.4
$ws1.Range("C122").Value = 40.62
$ws1.Range("D122").Value = -0.95

$ws1.Range("B123").Value = -1.39
$ws1.Range("C123").Value = 3.44
$ws1.Range("D123").Value = 21.55

$ws1.Range("B124").Value = 36.8
$ws1.Range("C124").Value = 38.27
$ws1.Range("D124").Value = -7.81

$ws1.Range("B125").Value = -36.86
$ws1.Range("C125").Value = -39.57
$ws1.Range("D125").Value = -12.81

$ws1.Range("B126").Value = 10.12
$ws1.Range("C126").Value = 7.19
$ws1.Range("D126").Value = 2.33
$ws1.Range("E126").Value = 37.04
$ws1.Range("F126").Value = 31.56
$ws1.Range("G126").Value = 0.16

# Add new row 127 (copy the style of the last dated row, then fill values)
$ws1.Range("A126").Copy()
$ws1.Range("A127").PasteSpecial(-4122)
$ws1.Range("A127").Value = 45809
$ws1.Range("B127").Value = -23.01
$ws1.Range("C127").Value = -26.28
$ws1.Range("D127").Value = 6.41
$ws1.Range("E127").Value = 12.16
$ws1.Range("F127").Value = 4.81
$ws1.Range("G127").Value = 5.6

# ---------------------------------------------------------------------------
# Sheet: "Medium Term"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Medium Term")

$ws2.Range("B107").Value = 21.36
$ws2.Range("C107").Value = 21.31
$ws2.Range("D107").Value = 13.94

$ws2.Range("B108").Value = 27.56
$ws2.Range("C108").Value = 25.12
$ws2.Range("D108").Value = 15.49

$ws2.Range("B109").Value = 37.43
$ws2.Range("C109").Value = 33.35
$ws2.Range("D109").Value = 20.61

$ws2.Range("C110").Value = 44.67
$ws2.Range("D110").Value = 29.85

$ws2.Range("C111").Value = 42.7
$ws2.Range("D111").Value = 28.35

$ws2.Range("B112").Value = 47.82
$ws2.Range("C112").Value = 43.02
$ws2.Range("D112").Value = 30.23

# Add new row 113 (copy the style of the last dated row, then fill values)
$ws2.Range("A112").Copy()
$ws2.Range("A113").PasteSpecial(-4122)
$ws2.Range("A113").Value = 45809
$ws2.Range("B113").Value = 16.54
$ws2.Range("C113").Value = 40.78
$ws2.Range("D113").Value = 31.68
